$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values
$headers = @("glycan", "binding_score", "monosaccharides", "motifs", "sasa", "flexibility", "has_multi_node_motifs")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows (glycan, binding_score, monosaccharides, motifs, sasa, flexibility, has_multi_node_motifs)
$data = @(
    @("Fuc(a1-2)Gal(b1-4)GlcNAc", 6.528127523312708, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 4.346787610201794, 1.057597373596624, $true),
    @("Fuc(a1-2)Gal(b1-4)[Fuc(a1-3)]GlcNAc", 4.162304469992177, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.822085034059461, 0.4456879832318703, $true),
    @("Fuc(a1-2)Gal(b1-4)[Fuc(a1-3)]GlcNAc(b1-3)[Fuc(a1-2)Gal(b1-4)[Fuc(a1-3)]GlcNAc(b1-6)]Gal(b1-4)Glc", 2.474660024377276, "['Fuc(a1-2)', 'Gal(b1-4)', 'Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 6.886293454185985, 14.90170621541335, $true),
    @("Fuc(a1-2)Gal(b1-4)[Fuc(a1-3)]GlcNAc(b1-3)[Fuc(a1-3)[Gal(b1-4)]GlcNAc(b1-6)]Gal(b1-4)Glc", 0.0828338083473565, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.727253196689237, 4.569905915454178, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)Glc", -0.07532320667243141, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.419172138929145, 1.022714362562547, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc", -0.37883078618528, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.496947085659928, 1.34727688863881, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)GalNAc", -0.197419301657273, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.346616244917194, 0.9918811636990688, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)[Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-6)]GalNAc", -0.2830121418291967, "['Fuc(a1-2)', 'Gal(b1-4)', 'Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 6.493081341348088, 7.342883445993724, $true),
    @("Fuc(a1-2)[GalNAc(a1-3)]Gal(b1-4)GlcNAc", -0.1750850165433314, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.362486045967387, 1.302091379100984, $true),
    @("Fuc(a1-2)[GalNAc(a1-3)]Gal(b1-4)GlcNAc(b1-3)GalNAc", -0.1914726002652398, "['Fuc(a1-2)', 'Gal(b1-4)']", "['Fuc(a1-2)Gal(b1-4)']", 3.375043197745113, 1.210960590183663, $true)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Build the bold + thin-border + center/top-aligned style on a single seed cell
# (keeps styles.xml minimal: 1 new font, 1 new border, 1 new cellXf), then
# propagate it via copy/paste-special-formats so no extra style entries leak.
$seed = $ws.Range("A1")
$seed.Font.Bold = $true
$seed.Borders.LineStyle = 1
$seed.HorizontalAlignment = -4108
$seed.VerticalAlignment = -4160

$seed.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$seed.Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

$excel.CutCopyMode = 0
